$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Row 1 originally has duplicated "value" header text in C1:F1 -- clear them.
$ws.Range("C1:F1").ClearContents()

# Remove the "Deletion" parameter row (was row 16: Deletion | 0 | 3).
$ws.Rows("16:16").Delete()

# Insert a new row after the (renamed) "production_function" row (row 8) for
# the new "L_curve" parameter.
$ws.Rows("9:9").Insert()

# Row 8 was "Model" -> now "production_function" (value stays "Sigmoid").
$ws.Range("A8").Value = "production_function"

# New row 9: L_curve | 1
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 1
$ws.Range("A9").Font.Name = "Verdana"
$ws.Range("A9").Font.Size = 10
$ws.Range("B9").NumberFormat = "0.00E+00"

# Make this sheet the active tab and move the selection to the new
# estimate_params value cell (B10, previously B9 before the insert).
$ws.Activate()
$ws.Range("B10").Select()
